$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "JSU(-0.8160485575245678, 1.0609090841005897, 0.7535327594610018, 2.7435104765436975)"
$ws.Range("C2").Value = "NIG(0.9182417664218405, 0.6388138035950721, 4.407928337252409, 5.05949476549884)"
$ws.Range("D2").Value = "MIE(7.446918536272837, 4.2549191529167505, -7.104222703372565, 8.936605002541208)"
$ws.Range("E2").Value = "NIG(1.6129452010259917, 1.302384818934952, 3.5659515727846722, 6.395753651967115)"
